# Actualización de la planilla de turnos
# Corrige los valores de "rut" que tenían espacios en blanco sobrantes
# al inicio/fin del texto, dejando el texto ya "limpio".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapa fila -> valor de rut corregido (sin espacios extra)
$rutFixes = @{
    4  = "16.366.615-4"
    7  = "16.759.697-5"
    10 = "15.343.687-8"
    11 = "16.751.516-9"
    22 = "15.343.687-8"
    42 = "15.343.687-8"
    50 = "16.366.615-4"
    55 = "16.790.074-7"
    73 = "15.343.687-8"
}

foreach ($row in $rutFixes.Keys) {
    $ws.Cells.Item($row, 2).Value = $rutFixes[$row]
}

# Actualiza la celda seleccionada de I11 a K11
$ws.Range("K11").Select()
